# Fixed Stimulus Absolute Timestamps
# Renames each task-order sheet (new timestamp suffixes) and updates the
# stimulus file names referenced in column B of each sheet.

$wb = $excel.ActiveWorkbook

# --- Rename sheets (tab names encode a run timestamp) ---------------------
$wb.Worksheets.Item(1).Name = "GNG_TO-1650477884350006"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778864404771"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477886446441"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778865054402"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778865664756"

# --- Sheet 1 (GNG_TO...) stimulus file names -------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778843220067.csv"
$ws1.Range("B3").Value = "GNG_stims-1650477884334004.csv"
$ws1.Range("B4").Value = "go_stims-16504778843350043.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778843490045.csv"

# --- Sheet 2 (NB_TO...) stimulus file names --------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650477885693008.csv"
$ws2.Range("B3").Value = "ZB-match_7-16504778847490065.csv"
$ws2.Range("B4").Value = "TB-1650477886214008.csv"
$ws2.Range("B5").Value = "TB-16504778864154425.csv"
$ws2.Range("B6").Value = "OB-16504778853370087.csv"
$ws2.Range("B7").Value = "OB-1650477885305012.csv"
$ws2.Range("B8").Value = "ZB-match_1-16504778845050418.csv"
$ws2.Range("B9").Value = "TB-16504778859310417.csv"
$ws2.Range("B10").Value = "ZB-match_1-16504778848970375.csv"

# --- Sheet 4 (TOL_TO...) stimulus file names -------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778864724753.csv"
$ws4.Range("B3").Value = "ZM_stims-1650477886449441.csv"
$ws4.Range("B4").Value = "MM_stims-1650477886488475.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778864734426.csv"
$ws4.Range("B6").Value = "MM_stims-16504778865044754.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477886489441.csv"

# --- Sheet 5 (vSAT_TO...) stimulus file names ------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504778865204728.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504778865514753.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778865354755.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778865084748.csv"
